$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the "MON Oct 30" timestamp used to be split across two runs
# (one run for "MON Oct 30", another for " 16:41:15 PDT 2017"). Collapse
# it back down to a single run by re-writing the found text in place.
# ---------------------------------------------------------------------
$mergeRange = $d.Content
$mergeFound = $mergeRange.Find.Execute(
    "MON Oct 30 16:41:15 PDT 2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "MON Oct 30 16:41:15 PDT 2017", 2)
if (-not $mergeFound) {
    throw "Could not find the 'MON Oct 30 16:41:15 PDT 2017' run pair to merge."
}

# ---------------------------------------------------------------------
# Change 2: a new "TUE OCT 31" purchase-details entry (MAMATHA / HN T
# chick-in, 2/11/2017) needs to be appended right after the existing
# "- ACC 27/10/2017" line.
#
# NOTE: this runtime's PowerShell-alike only binds function parameters
# positionally, so every helper below is called positionally (no
# "-Name value" switches/parameters).
# ---------------------------------------------------------------------

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-CourierRPr($Bold) {
    $boldTag = ""
    if ($Bold) { $boldTag = "<w:b/>" }
    return "<w:rPr><w:rFonts w:ascii=`"Courier New`" w:hAnsi=`"Courier New`" w:cs=`"Courier New`"/>$boldTag</w:rPr>"
}

function New-PlainTextParagraphXml($RunsXml, $Bold) {
    $pPrRPr = Get-CourierRPr $Bold
    return "<w:p $wNs><w:pPr><w:pStyle w:val=`"PlainText`"/>$pPrRPr</w:pPr>$RunsXml</w:p>"
}

function New-TextRunXml($Text, $Bold, $PreserveSpace) {
    $rPr = Get-CourierRPr $Bold
    $spaceAttr = ""
    if ($PreserveSpace) { $spaceAttr = ' xml:space="preserve"' }
    $escaped = $Text -replace '&', '&amp;'
    return "<w:r>$rPr<w:t$spaceAttr>$escaped</w:t></w:r>"
}

function New-TabRunXml($Bold) {
    $rPr = Get-CourierRPr $Bold
    return "<w:r>$rPr<w:tab/></w:r>"
}

function New-TabTextRunXml($Text, $Bold) {
    $rPr = Get-CourierRPr $Bold
    $escaped = $Text -replace '&', '&amp;'
    return "<w:r>$rPr<w:tab/><w:t>$escaped</w:t></w:r>"
}

# A "label" line: the label text, a run of plain tabs, then one final
# run that bundles the last tab together with the "- value" text -
# exactly how Word splits runs when the line is typed with trailing tabs.
function New-LabelValueParagraphXml($Label, $TabCount, $Value, $Bold) {
    $runs = New-TextRunXml $Label $Bold $false
    for ($i = 0; $i -lt $TabCount; $i++) {
        $runs += New-TabRunXml $Bold
    }
    $runs += New-TabTextRunXml $Value $Bold
    return New-PlainTextParagraphXml $runs $Bold
}

$rows = @(
    @{ Label = "Person Name";       Tabs = 3; Value = "- HN T" },
    @{ Label = "Bill number";       Tabs = 3; Value = "- 517" },
    @{ Separator = $true },
    @{ Label = "Item Name";         Tabs = 3; Value = "- CARROT" },
    @{ Label = "Number of Pockets"; Tabs = 2; Value = "- 1" },
    @{ Label = "Number of KGs";     Tabs = 2; Value = "- 91" },
    @{ Label = "Rate";              Tabs = 4; Value = "- 52" },
    @{ Label = "Transport & Miscellaneous"; Tabs = 0; Value = "- 115" },
    @{ Label = "Total Price";       Tabs = 3; Value = "- 4847.0" },
    @{ Label = "Amount balance";    Tabs = 2; Value = "- 10813.0"; Bold = $true }
)

$separatorLine = "---------------------------------------------------------------"

$blockXml = ""

# Blank separator paragraph before the new entry.
$blockXml += New-PlainTextParagraphXml "" $false

# Timestamp line, split into two runs just like the other date/time
# stamps already present in the document.
$timestampRuns = (New-TextRunXml "TUE OCT 31" $false $false) + (New-TextRunXml " 17:10:21 PDT 2017" $false $true)
$blockXml += New-PlainTextParagraphXml $timestampRuns $false

foreach ($row in $rows) {
    if ($row.Separator) {
        $blockXml += New-PlainTextParagraphXml (New-TextRunXml $separatorLine $false $false) $false
        continue
    }
    $blockXml += New-LabelValueParagraphXml $row.Label $row.Tabs $row.Value ([bool]$row.Bold)
}

# Trailing blank paragraphs (bold run-properties, then plain) that close
# out the new block, matching the spacing used elsewhere in the file.
$blockXml += New-PlainTextParagraphXml "" $true
$blockXml += New-PlainTextParagraphXml "" $false

# Anchor on the existing "- ACC 27/10/2017" line and insert the new
# block immediately after it (collapsed insertion point so the existing
# paragraph's own text/content is left completely untouched).
$anchor = $d.Content
$anchorFound = $anchor.Find.Execute(
    "- ACC 27/10/2017", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $anchorFound) {
    throw "Could not find the '- ACC 27/10/2017' anchor line."
}
$insertionPoint = $d.Range($anchor.End, $anchor.End)
[void]$insertionPoint.InsertXML($blockXml)
